$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (D and E), shifting existing
# quarterly data (old D:K) to the right (new F:M), to make room for the
# two newest fiscal-quarter columns of financial data.
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats / styles from the (now shifted) old column D --
# which now lives in column F -- into the two freshly inserted columns
# so the new cells inherit the correct date / number formatting instead
# of defaulting to the General style copied from column C.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the latest two quarters' financial data.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 639500
$ws.Cells.Item(8, 5).Value = 647300
$ws.Cells.Item(9, 4).Value = 481300
$ws.Cells.Item(9, 5).Value = 485400
$ws.Cells.Item(10, 4).Value = 158200
$ws.Cells.Item(10, 5).Value = 161900
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 574400
$ws.Cells.Item(17, 5).Value = 581100
$ws.Cells.Item(18, 4).Value = 65100
$ws.Cells.Item(18, 5).Value = 66200
$ws.Cells.Item(20, 4).Value = 300
$ws.Cells.Item(20, 5).Value = 200
$ws.Cells.Item(21, 4).Value = 77700
$ws.Cells.Item(21, 5).Value = 78300
$ws.Cells.Item(22, 4).Value = 9700
$ws.Cells.Item(22, 5).Value = 9400
$ws.Cells.Item(23, 4).Value = 55800
$ws.Cells.Item(23, 5).Value = 57000
$ws.Cells.Item(24, 4).Value = 17200
$ws.Cells.Item(24, 5).Value = 14400
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 38600
$ws.Cells.Item(26, 5).Value = 42700
$ws.Cells.Item(27, 4).Value = 38600
$ws.Cells.Item(27, 5).Value = 42700
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -300
$ws.Cells.Item(32, 5).Value = -200
$ws.Cells.Item(33, 4).Value = 38600
$ws.Cells.Item(33, 5).Value = 42700
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 38600
$ws.Cells.Item(35, 5).Value = 42700
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 100900
$ws.Cells.Item(41, 5).Value = 93500
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 407100
$ws.Cells.Item(43, 5).Value = 419700
$ws.Cells.Item(44, 4).Value = 169000
$ws.Cells.Item(44, 5).Value = 161900
$ws.Cells.Item(45, 4).Value = 27700
$ws.Cells.Item(45, 5).Value = 24100
$ws.Cells.Item(46, 4).Value = 704700
$ws.Cells.Item(46, 5).Value = 699100
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 168000
$ws.Cells.Item(48, 5).Value = 166700
$ws.Cells.Item(49, 4).Value = 1563400
$ws.Cells.Item(49, 5).Value = 1567900
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 18500
$ws.Cells.Item(52, 5).Value = 23100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 2454500
$ws.Cells.Item(54, 5).Value = 2456800
$ws.Cells.Item(57, 4).Value = 313200
$ws.Cells.Item(57, 5).Value = 300900
$ws.Cells.Item(58, 4).Value = 26900
$ws.Cells.Item(58, 5).Value = 23400
$ws.Cells.Item(59, 4).Value = 104200
$ws.Cells.Item(59, 5).Value = 116200
$ws.Cells.Item(60, 4).Value = 444300
$ws.Cells.Item(60, 5).Value = 440600
$ws.Cells.Item(61, 4).Value = 716600
$ws.Cells.Item(61, 5).Value = 718400
$ws.Cells.Item(62, 4).Value = 221600
$ws.Cells.Item(62, 5).Value = 211700
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 1382400
$ws.Cells.Item(66, 5).Value = 1370800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 441900
$ws.Cells.Item(72, 5).Value = 403300
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 1072100
$ws.Cells.Item(76, 5).Value = 1086000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 38600
$ws.Cells.Item(81, 5).Value = 42700
$ws.Cells.Item(83, 4).Value = 12300
$ws.Cells.Item(83, 5).Value = 11900
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 71100
$ws.Cells.Item(89, 5).Value = 54600
$ws.Cells.Item(91, 4).Value = -10100
$ws.Cells.Item(91, 5).Value = -14900
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -9300
$ws.Cells.Item(94, 5).Value = -16400
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -54400
$ws.Cells.Item(100, 5).Value = -10500
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 7500
$ws.Cells.Item(102, 5).Value = 27700

